$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.376.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "'1.844.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.15%  "
$ws.Range("D5").Value = "'315.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("D7").Value = "'0.4724"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("D8").Value = "'0.3702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.07454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.04%  "
$ws.Range("D10").Value = "'0.8856"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'20.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'1.826.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "'0.07352"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.98%  "
$ws.Range("D14").Value = "'5.492"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "'93.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "'6.597"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "'0.000008847"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "'14.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "'27.400.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").Value = "'5.332"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'10.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'2.083.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "'1.906"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "'153.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").Value = "'18.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "'2.190"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'118.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").Value = "'0.08966"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'0.7620"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").Value = "'1.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "'4.575"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "'2.940"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "'0.05359"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "'0.01964"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'3.008"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").Value = "'7.366"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'2.403"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "'0.5371"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.1669"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'8.597"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("D46").Value = "'0.4980"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'10.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D49").Value = "'1.685"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'104.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").Value = "'0.06327"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
